$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that currently sits between
#    "...hvis det ikke er muligt " and "– upload " in the first bullet under
#    "Overordnet ansvar". It will be re-created later at its new location
#    (inside the newly inserted bullet about "Like"-ing Facebook posts).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert a brand-new bulleted paragraph right before the
#    "Vi skal være gode til at give og modtage konstruktiv kritik." bullet,
#    with the new sentence about remembering to "Like" Facebook posts.
# ---------------------------------------------------------------------------
$critiqueRange = $d.Content
$critiqueRange.Find.Execute("Vi skal være gode til at give og modtage konstruktiv kritik.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$critiqueIndex = $critiqueRange.Paragraphs(1).Index

$critiqueParagraph = $d.Paragraphs($critiqueIndex)
$critiqueParagraph.Range.InsertParagraphBefore()

# After InsertParagraphBefore, the new (empty) paragraph takes over the
# original index; the old content shifted one slot down.
$newParagraph = $d.Paragraphs($critiqueIndex)
$newParagraph.Range.Text = "Man skal huske at ”Like” opsalg på Facebook når man har læst dem."

# Place the "_GoBack" bookmark right after "Facebook" in the new sentence.
$bmFind = $newParagraph.Range
$bmFind.Find.Execute("Facebook", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$bmPoint = $d.Range($bmFind.End, $bmFind.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)
